$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = "29.104.36"
$d.Style = "Normal"
$ws.Range("E2").Value = "  +1.05%  "

$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = "1.902.17"
$d.Style = "Normal"
$ws.Range("E3").Value = "  +1.15%  "

$d = $ws.Range("D4")
$d.NumberFormat = "@"
$d.Value = "1.003"
$d.Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = "327.93"
$d.Style = "Normal"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("E6").Value = "  -0.08%  "

$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = "0.4655"
$d.Style = "Normal"
$ws.Range("E7").Value = "  -0.46%  "

$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = "0.3928"
$d.Style = "Normal"
$ws.Range("E8").Value = "  -0.16%  "

$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = "47.09"
$d.Style = "Normal"
$ws.Range("E9").Value = "  +1.23%  "

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = "0.07966"
$d.Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = "1.014"
$d.Style = "Normal"

$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = "22.17"
$d.Style = "Normal"
$ws.Range("E12").Value = "  -0.73%  "

$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = "1.931.49"
$d.Style = "Normal"
$ws.Range("E13").Value = "  +4.62%  "

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "7.150"
$d.Style = "Normal"
$ws.Range("E14").Value = "  +1.91%  "

$ws.Range("E15").Value = "  +0.63%  "

$ws.Range("E16").Value = "  +0.30%  "

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = "89.58"
$d.Style = "Normal"
$ws.Range("E17").Value = "  +0.96%  "

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = "1.004"
$d.Style = "Normal"
$ws.Range("E18").Value = "  -0.20%  "

$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = "0.00001015"
$d.Style = "Normal"
$ws.Range("E19").Value = "  +0.52%  "

$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = "17.29"
$d.Style = "Normal"
$ws.Range("E20").Value = "  +1.86%  "

$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = "1.003"
$d.Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = "29.076.47"
$d.Style = "Normal"
$ws.Range("E22").Value = "  +0.94%  "

$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = "5.361"
$d.Style = "Normal"
$ws.Range("E23").Value = "  +0.31%  "

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = "11.15"
$d.Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = "2.146.57"
$d.Style = "Normal"
$ws.Range("E25").Value = "  +3.24%  "

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = "2.065"
$d.Style = "Normal"
$ws.Range("E26").Value = "  -2.51%  "

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = "155.38"
$d.Style = "Normal"
$ws.Range("E27").Value = "  +0.99%  "

$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = "19.85"
$d.Style = "Normal"
$ws.Range("E28").Value = "  +2.29%  "

$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = "5.898"
$d.Style = "Normal"
$ws.Range("E29").Value = "  +2.59%  "

$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = "1.995"
$d.Style = "Normal"
$ws.Range("E30").Value = "  -0.35%  "

$d = $ws.Range("D31")
$d.NumberFormat = "@"
$d.Value = "120.47"
$d.Style = "Normal"
$ws.Range("E31").Value = "  +0.42%  "

$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = "0.09398"
$d.Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "

$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = "0.9434"
$d.Style = "Normal"
$ws.Range("E33").Value = "  +0.42%  "

$ws.Range("E34").Value = "  +0.87%  "

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = "1.357"
$d.Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "

$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = "3.261"
$d.Style = "Normal"
$ws.Range("E36").Value = "  -2.74%  "

$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = "0.05852"
$d.Style = "Normal"
$ws.Range("E37").Value = "  -1.12%  "

$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = "1.181"
$d.Style = "Normal"
$ws.Range("E38").Value = "  +1.96%  "

$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = "8.111"
$d.Style = "Normal"
$ws.Range("E39").Value = "  +2.82%  "

$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = "0.02103"
$d.Style = "Normal"
$ws.Range("E40").Value = "  -0.90%  "

$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = "0.5846"
$d.Style = "Normal"
$ws.Range("E41").Value = "  +2.08%  "

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = "1.004"
$d.Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = "0.1819"
$d.Style = "Normal"
$ws.Range("E43").Value = "  +1.25%  "

$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = "10.03"
$d.Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "

$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = "2.293"
$d.Style = "Normal"
$ws.Range("E45").Value = "  +8.34%  "

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = "0.5477"
$d.Style = "Normal"
$ws.Range("E46").Value = "  +2.58%  "

$ws.Range("E47").Value = "  +0.66%  "

$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = "0.07215"
$d.Style = "Normal"
$ws.Range("E48").Value = "  -1.53%  "

$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = "1.880"
$d.Style = "Normal"
$ws.Range("E49").Value = "  +1.79%  "

$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = "1.120"
$d.Style = "Normal"
$ws.Range("E50").Value = "  -3.03%  "

$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = "113.39"
$d.Style = "Normal"
$ws.Range("E51").Value = "  -0.63%  "
